$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 3
Write-Host "ScrollRow: $($excel.ActiveWindow.ScrollRow)"
Write-Host "ScrollColumn: $($excel.ActiveWindow.ScrollColumn)"
